$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 814
$ws.Range("F5").Value = 499
$ws.Range("F7").Value = 477
$ws.Range("F8").Value = 887
$ws.Range("F9").Value = 116
$ws.Range("F10").Value = 820
$ws.Range("F11").Value = 645
$ws.Range("F12").Value = 119
$ws.Range("F14").Value = 52
$ws.Range("F16").Value = 214
$ws.Range("F18").Value = 463
$ws.Range("F19").Value = 1229
$ws.Range("F21").Value = 941
$ws.Range("F22").Value = 2687
$ws.Range("F23").Value = 1156
$ws.Range("F25").Value = 140
$ws.Range("F26").Value = 1198
$ws.Range("F28").Value = 894
$ws.Range("F29").Value = 94
$ws.Range("F30").Value = 1228

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 481

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 814
$ws.Range("F6").Value = 499
$ws.Range("F9").Value = 477
$ws.Range("F10").Value = 481
$ws.Range("F14").Value = 887
$ws.Range("F15").Value = 116
$ws.Range("F17").Value = 645
$ws.Range("F18").Value = 119
$ws.Range("F24").Value = 52
$ws.Range("F27").Value = 214
$ws.Range("F29").Value = 463
$ws.Range("F30").Value = 1229
$ws.Range("F32").Value = 941
$ws.Range("F33").Value = 2687
$ws.Range("F34").Value = 1156
$ws.Range("F36").Value = 140
$ws.Range("F37").Value = 1198
$ws.Range("F40").Value = 894
$ws.Range("F41").Value = 94
$ws.Range("F42").Value = 1228

